# Leave Card update — insert a new leave-credit row (UT 0-0-8 / UT 0-0-2)
# above the old row 44 ("2023" header row) on the "2018 LEAVE CREDITS" sheet.
#
# Net effect (see commit diff):
#   - A brand new physical row is inserted at row 44, pushing the former
#     rows 44-135 down to 45-136.
#   - The (new) row 43 becomes a small "Absence Undertime W/ Pay" entry
#     (UT(0-0-2), 0.004) with no PERIOD/EARNED value.
#   - The (new) row 44 becomes a dated EARNED entry (UT0-0-8), 1.25 earned,
#     0.017 absence undertime w/ pay).
#   - Table1 grows from A8:K135 to A8:K136.
#   - CONVERTION!F3 changes from 14 to 2 (recalculates G3 automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws.Activate()

# --- 1. Insert a new physical row at 44 (shifts old 44..135 -> 45..136) ---
$ws.Rows.Item(44).Insert()

# Copy the formatting from row 43 (the row immediately above the insertion
# point) into the freshly inserted row 44 so styles match the rest of the
# table (only the used A:K range, not the whole 16384-column row).
$ws.Range("A43:K43").Copy()
$ws.Range("A44:K44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Populate the new row 44 as the "UT0-0-8)" earned entry ---
# (set before row 43's string so the new shared-string table gets the two
# new entries in the same order the original author's file has them: index
# 74 = "UT0-0-8)", index 75 = "UT(0-0-2)")
$ws.Range("A44").Value = 44926
$ws.Range("B44").Value = "UT0-0-8)"
$ws.Range("C44").Value = 1.25
$ws.Range("D44").Value = 0.017
$ws.Range("G44").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 3. Re-purpose row 43 as the "UT(0-0-2)" absence-undertime entry ---
$ws.Range("A43").ClearContents()
$ws.Range("B43").Value = "UT(0-0-2)"
$ws.Range("C43").ClearContents()
$ws.Range("D43").Value = 0.004
$ws.Range("G43").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 4. Grow Table1 so it once again covers the whole data range ---
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K136"))

# Row 136 (the new trailing blank/total-style row pushed out of the old
# table bounds) needs its calculated-column formula re-applied explicitly;
# the engine leaves it as a broken non-table reference after the resize.
$ws.Range("G136").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 5. CONVERTION sheet: DAYS lookup input changes from 14 to 2 ---
$cv = $wb.Worksheets.Item("CONVERTION")
$cv.Range("F3").Value = 2

# --- 6. Recalculate everything so all cached formula results are fresh ---
$excel.CalculateFull()

# --- 7. Leave the cursor where the author left it ---
$ws.Range("F46").Select()
